$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in H1, matching the bold/centered header style
# used by the rest of row 1 (same style as G1).
$ws.Range("H1").Value = "Save"
$ws.Range("H1").Style = $ws.Range("G1").Style

# Fill in the new "Save" data column (H2:H9)
$saveValues = @(0, 1, 0, 0, 0, 0, 0, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
